$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced to
# text format first, otherwise Excel auto-converts them to numeric values
# (e.g. "1.00" -> 1) and loses the original formatted-text representation.
$ws.Range("D2").Value = "64.174.52"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.322.17"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.83"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.71"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.313.27"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  +8.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.629"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.59"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("E13").Value = "  +4.41%  "
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "3.855.14"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.331.43"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "65.109.61"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.69"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.983"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "453.29"
$ws.Range("E22").Value = "  +7.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.95"
$ws.Range("E23").Value = "  +6.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.08"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.74"
$ws.Range("E25").Value = "  +3.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.73"
$ws.Range("E26").Value = "  +6.74%  "
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.68"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("E29").Value = "  +5.63%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.54"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.42"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "566.31"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "60.56"
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.56"
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.18"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "0.0₃0733"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "3.058.43"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.76"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.64"
$ws.Range("E49").Value = "  +6.72%  "
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.14"
$ws.Range("E51").Value = "  +0.70%  "
